# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 28978
$ws.Range("J3").Value = 28978
$ws.Range("L3").Value = 28978
$ws.Range("N3").Value = -29206
$ws.Range("H28").Value = 154.09525
$ws.Range("J28").Value = 169
$ws.Range("L28").Value = 169
$ws.Range("N28").Value = -1139
$ws.Range("H40").Value = 2199.6191
$ws.Range("I40").Value = 1622.4615
$ws.Range("J40").Value = 3137.5
$ws.Range("K40").Value = 1622.4615
$ws.Range("L40").Value = 3137.5
$ws.Range("M40").Value = -1447.4615
$ws.Range("N40").Value = -3487.5
$ws.Range("H64").Value = 3223.5
$ws.Range("I64").Value = 3222
$ws.Range("J64").Value = 3226.125
$ws.Range("K64").Value = 3222
$ws.Range("L64").Value = 3226.125
$ws.Range("M64").Value = -2974
$ws.Range("N64").Value = -3722.125
$ws.Range("H67").Value = 3223.5
$ws.Range("I67").Value = 3222
$ws.Range("J67").Value = 3226.125
$ws.Range("K67").Value = 3222
$ws.Range("L67").Value = 3226.125
$ws.Range("M67").Value = -2364
$ws.Range("N67").Value = -4942.125
$ws.Range("H76").Value = 3499.158
$ws.Range("I76").Value = 3211.111
$ws.Range("J76").Value = 3758.4
$ws.Range("K76").Value = 3211.111
$ws.Range("L76").Value = 3758.4
$ws.Range("M76").Value = -2896.111
$ws.Range("N76").Value = -4388.4
$ws.Range("H79").Value = 3499.158
$ws.Range("I79").Value = 3211.111
$ws.Range("J79").Value = 3758.4
$ws.Range("K79").Value = 3211.111
$ws.Range("L79").Value = 3758.4
$ws.Range("M79").Value = -2119.111
$ws.Range("N79").Value = -5942.4
$ws.Range("H98").Value = 917.5
$ws.Range("I98").Value = 988.6875
$ws.Range("J98").Value = 348
$ws.Range("K98").Value = 988.6875
$ws.Range("L98").Value = 348
$ws.Range("M98").Value = 509.3125
$ws.Range("N98").Value = -3344
$ws.Range("H102").Value = 28978
$ws.Range("J102").Value = 28978
$ws.Range("L102").Value = 28978
$ws.Range("N102").Value = -35468
$ws.Range("H122").Value = 917.5
$ws.Range("I122").Value = 988.6875
$ws.Range("J122").Value = 348
$ws.Range("K122").Value = 2966.0625
$ws.Range("L122").Value = 1044
$ws.Range("M122").Value = -516.0625
$ws.Range("N122").Value = -5944
$ws.Range("H137").Value = 27307.236
$ws.Range("I137").Value = 786.5217
$ws.Range("J137").Value = 67972.336
$ws.Range("K137").Value = 2359.5651
$ws.Range("L137").Value = 203917.008
$ws.Range("M137").Value = 190.4349000000002
$ws.Range("N137").Value = -209017.008
$ws.Range("H138").Value = 4020.6416
$ws.Range("I138").Value = 2971
$ws.Range("J138").Value = 4397.436
$ws.Range("K138").Value = 8913
$ws.Range("L138").Value = 13192.308
$ws.Range("M138").Value = -3773
$ws.Range("N138").Value = -23472.308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1973.0769
$ws.Range("I61").Value = 1676.4706
$ws.Range("J61").Value = 2533.3333
$ws.Range("K61").Value = 1676.4706
$ws.Range("L61").Value = 2533.3333
$ws.Range("M61").Value = -1464.4706
$ws.Range("N61").Value = -2957.3333
$ws.Range("H63").Value = 2161
$ws.Range("I63").Value = 2161
$ws.Range("K63").Value = 2161
$ws.Range("M63").Value = -1475
$ws.Range("H66").Value = 2161
$ws.Range("I66").Value = 2161
$ws.Range("K66").Value = 10805
$ws.Range("M66").Value = -7373
$ws.Range("H88").Value = 200161090
$ws.Range("I88").Value = 1200
$ws.Range("J88").Value = 250201060
$ws.Range("K88").Value = 1200
$ws.Range("L88").Value = 250201060
$ws.Range("M88").Value = -794
$ws.Range("N88").Value = -250201872
$ws.Range("H91").Value = 200161090
$ws.Range("I91").Value = 1200
$ws.Range("J91").Value = 250201060
$ws.Range("K91").Value = 1200
$ws.Range("L91").Value = 250201060
$ws.Range("M91").Value = 204
$ws.Range("N91").Value = -250203868
$ws.Range("H92").Value = 34000
$ws.Range("J92").Value = 34000
$ws.Range("L92").Value = 34000
$ws.Range("N92").Value = -38992
$ws.Range("H96").Value = 22896
$ws.Range("J96").Value = 22896
$ws.Range("L96").Value = 22896
$ws.Range("N96").Value = -28388
$ws.Range("H110").Value = 1633.44
$ws.Range("I110").Value = 1554.619
$ws.Range("J110").Value = 2047.25
$ws.Range("K110").Value = 1554.619
$ws.Range("L110").Value = 2047.25
$ws.Range("M110").Value = 490.3810000000001
$ws.Range("N110").Value = -6137.25
$ws.Range("H132").Value = 33272.312
$ws.Range("I132").Value = 45144
$ws.Range("J132").Value = 2933.5557
$ws.Range("K132").Value = 135432
$ws.Range("L132").Value = 8800.667099999999
$ws.Range("M132").Value = -132902
$ws.Range("N132").Value = -13860.6671
$ws.Range("H136").Value = 1973.0769
$ws.Range("I136").Value = 1676.4706
$ws.Range("J136").Value = 2533.3333
$ws.Range("K136").Value = 5029.4118
$ws.Range("L136").Value = 7599.999899999999
$ws.Range("M136").Value = -2479.4118
$ws.Range("N136").Value = -12699.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 815.5
$ws.Range("I24").Value = 595.625
$ws.Range("J24").Value = 1695
$ws.Range("K24").Value = 595.625
$ws.Range("L24").Value = 1695
$ws.Range("M24").Value = -360.625
$ws.Range("N24").Value = -2165
$ws.Range("H94").Value = 6327.6113
$ws.Range("I94").Value = 756.7273
$ws.Range("K94").Value = 756.7273
$ws.Range("M94").Value = -305.7273
$ws.Range("H105").Value = 1812.5555
$ws.Range("I105").Value = 1527.5883
$ws.Range("J105").Value = 2297
$ws.Range("K105").Value = 1527.5883
$ws.Range("L105").Value = 2297
$ws.Range("M105").Value = 219.4117000000001
$ws.Range("N105").Value = -5791

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 15966.667
$ws.Range("J43").Value = 15966.667
$ws.Range("L43").Value = 15966.667
$ws.Range("N43").Value = -16334.667
$ws.Range("H62").Value = 3090.647
$ws.Range("J62").Value = 3286.3333
$ws.Range("L62").Value = 3286.3333
$ws.Range("N62").Value = -4534.3333
$ws.Range("H65").Value = 3090.647
$ws.Range("J65").Value = 3286.3333
$ws.Range("L65").Value = 16431.6665
$ws.Range("N65").Value = -22671.6665
$ws.Range("H101").Value = 15966.667
$ws.Range("J101").Value = 15966.667
$ws.Range("L101").Value = 15966.667
$ws.Range("N101").Value = -22456.667
$ws.Range("H132").Value = 2805.6843
$ws.Range("I132").Value = 5670.3335
$ws.Range("J132").Value = 2268.5625
$ws.Range("K132").Value = 17011.0005
$ws.Range("L132").Value = 6805.6875
$ws.Range("M132").Value = -14481.0005
$ws.Range("N132").Value = -11865.6875
$ws.Range("H134").Value = 2454.4707
$ws.Range("I134").Value = 2154.1333
$ws.Range("J134").Value = 4707
$ws.Range("K134").Value = 6462.3999
$ws.Range("L134").Value = 14121
$ws.Range("M134").Value = -3927.3999
$ws.Range("N134").Value = -19191

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1172.9
$ws.Range("I136").Value = 1132.25
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 3396.75
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = 1703.25
$ws.Range("N136").Value = -13800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 92161864
$ws.Range("I70").Value = 138238620
$ws.Range("J70").Value = 8333.333000000001
$ws.Range("K70").Value = 138238620
$ws.Range("L70").Value = 8333.333000000001
$ws.Range("M70").Value = -138238350
$ws.Range("N70").Value = -8873.333000000001
$ws.Range("H73").Value = 92161864
$ws.Range("I73").Value = 138238620
$ws.Range("J73").Value = 8333.333000000001
$ws.Range("K73").Value = 138238620
$ws.Range("L73").Value = 8333.333000000001
$ws.Range("M73").Value = -138237684
$ws.Range("N73").Value = -10205.333
$ws.Range("H80").Value = 3756.5908
$ws.Range("I80").Value = 4666.5
$ws.Range("J80").Value = 2998.3333
$ws.Range("K80").Value = 4666.5
$ws.Range("L80").Value = 2998.3333
$ws.Range("M80").Value = -3668.5
$ws.Range("N80").Value = -4994.3333
$ws.Range("H83").Value = 3756.5908
$ws.Range("I83").Value = 4666.5
$ws.Range("J83").Value = 2998.3333
$ws.Range("K83").Value = 23332.5
$ws.Range("L83").Value = 14991.6665
$ws.Range("M83").Value = -18340.5
$ws.Range("N83").Value = -24975.6665
$ws.Range("H101").Value = 38000
$ws.Range("J101").Value = 38000
$ws.Range("L101").Value = 38000
$ws.Range("N101").Value = -44490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2223065.5
$ws.Range("I22").Value = 3333773
$ws.Range("J22").Value = 1650
$ws.Range("K22").Value = 3333773
$ws.Range("L22").Value = 1650
$ws.Range("M22").Value = -3333478
$ws.Range("N22").Value = -2240
$ws.Range("H27").Value = 2223065.5
$ws.Range("I27").Value = 3333773
$ws.Range("J27").Value = 1650
$ws.Range("K27").Value = 3333773
$ws.Range("L27").Value = 1650
$ws.Range("M27").Value = -3333666
$ws.Range("N27").Value = -1864
$ws.Range("H40").Value = 1717
$ws.Range("I40").Value = 1346.25
$ws.Range("J40").Value = 3200
$ws.Range("K40").Value = 1346.25
$ws.Range("L40").Value = 3200
$ws.Range("M40").Value = -1210.25
$ws.Range("N40").Value = -3472
$ws.Range("H58").Value = 6733.3335
$ws.Range("I58").Value = 4200
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 4200
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = -3940
$ws.Range("N58").Value = -8520
$ws.Range("H93").Value = 1429.381
$ws.Range("I93").Value = 1433.4615
$ws.Range("J93").Value = 1422.75
$ws.Range("K93").Value = 1433.4615
$ws.Range("L93").Value = 1422.75
$ws.Range("M93").Value = -185.4614999999999
$ws.Range("N93").Value = -3918.75
$ws.Range("H122").Value = 4108.6855
$ws.Range("I122").Value = 4362.143
$ws.Range("J122").Value = 3728.5
$ws.Range("K122").Value = 13086.429
$ws.Range("L122").Value = 11185.5
$ws.Range("M122").Value = -10636.429
$ws.Range("N122").Value = -16085.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4995.3335
$ws.Range("J2").Value = 4995.3335
$ws.Range("L2").Value = 4995.3335
$ws.Range("N2").Value = -5219.3335
$ws.Range("H5").Value = 10000000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H18").Value = 1500
$ws.Range("J18").Value = 1500
$ws.Range("L18").Value = 1500
$ws.Range("N18").Value = -1846
$ws.Range("H107").Value = 2121.9167
$ws.Range("I107").Value = 1732.5
$ws.Range("J107").Value = 2900.75
$ws.Range("K107").Value = 5197.5
$ws.Range("L107").Value = 8702.25
$ws.Range("M107").Value = -3277.5
$ws.Range("N107").Value = -12542.25
